$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "green transition"
$ws.Range("C2").Value = "greenhouse effect"
$ws.Range("B3").Value = "green transition"
$ws.Range("C3").Value = "loss of biodiversity"
$ws.Range("B4").Value = "green transition"
$ws.Range("C4").Value = "extreme weather events"
$ws.Range("B5").Value = "green transition"
$ws.Range("D5").Value = 2
$ws.Range("B6").Value = "green transition"
$ws.Range("C6").Value = "emissions"
$ws.Range("D6").Value = 7
$ws.Range("B7").Value = "green transition"
$ws.Range("C7").Value = "global warming"
$ws.Range("B8").Value = "green transition"
$ws.Range("C8").Value = "melting glaciers"
$ws.Range("D8").Value = 0
$ws.Range("B9").Value = "green transition"
$ws.Range("C9").Value = "renewable energy"
$ws.Range("D9").Value = 6
$ws.Range("B10").Value = "green transition"
$ws.Range("C10").Value = "misinformation"
$ws.Range("B11").Value = "greenhouse effect"
$ws.Range("C11").Value = "loss of biodiversity"
$ws.Range("D11").Value = 0
$ws.Range("B12").Value = "greenhouse effect"
$ws.Range("C12").Value = "extreme weather events"
$ws.Range("B13").Value = "greenhouse effect"
$ws.Range("D13").Value = 80
$ws.Range("B14").Value = "greenhouse effect"
$ws.Range("C14").Value = "emissions"
$ws.Range("B15").Value = "greenhouse effect"
$ws.Range("C15").Value = "global warming"
$ws.Range("B16").Value = "greenhouse effect"
$ws.Range("C16").Value = "melting glaciers"
$ws.Range("D16").Value = 0
$ws.Range("B17").Value = "greenhouse effect"
$ws.Range("C17").Value = "renewable energy"
$ws.Range("B18").Value = "greenhouse effect"
$ws.Range("C18").Value = "misinformation"
$ws.Range("B19").Value = "loss of biodiversity"
$ws.Range("C19").Value = "extreme weather events"
$ws.Range("D19").Value = 1
$ws.Range("B20").Value = "loss of biodiversity"
$ws.Range("D20").Value = 1
$ws.Range("B21").Value = "loss of biodiversity"
$ws.Range("C21").Value = "emissions"
$ws.Range("D21").Value = 7
$ws.Range("B22").Value = "loss of biodiversity"
$ws.Range("C22").Value = "global warming"
$ws.Range("D22").Value = 2
$ws.Range("B23").Value = "loss of biodiversity"
$ws.Range("C23").Value = "melting glaciers"
$ws.Range("B24").Value = "loss of biodiversity"
$ws.Range("C24").Value = "renewable energy"
$ws.Range("B25").Value = "loss of biodiversity"
$ws.Range("C25").Value = "misinformation"
$ws.Range("B26").Value = "extreme weather events"
$ws.Range("B27").Value = "extreme weather events"
$ws.Range("C27").Value = "emissions"
$ws.Range("D27").Value = 16
$ws.Range("B28").Value = "extreme weather events"
$ws.Range("C28").Value = "global warming"
$ws.Range("B29").Value = "extreme weather events"
$ws.Range("C29").Value = "melting glaciers"
$ws.Range("D29").Value = 4
$ws.Range("B30").Value = "extreme weather events"
$ws.Range("C30").Value = "renewable energy"
$ws.Range("B31").Value = "extreme weather events"
$ws.Range("C31").Value = "misinformation"
$ws.Range("D31").Value = 4
$ws.Range("C32").Value = "emissions"
$ws.Range("D32").Value = 12
$ws.Range("C33").Value = "global warming"
$ws.Range("C34").Value = "melting glaciers"
$ws.Range("D34").Value = 0
$ws.Range("C35").Value = "renewable energy"
$ws.Range("C36").Value = "misinformation"
$ws.Range("B37").Value = "emissions"
$ws.Range("C37").Value = "global warming"
$ws.Range("D37").Value = 19
$ws.Range("B38").Value = "emissions"
$ws.Range("C38").Value = "melting glaciers"
$ws.Range("D38").Value = 6
$ws.Range("B39").Value = "emissions"
$ws.Range("C39").Value = "renewable energy"
$ws.Range("B40").Value = "emissions"
$ws.Range("C40").Value = "misinformation"
$ws.Range("B41").Value = "global warming"
$ws.Range("C41").Value = "melting glaciers"
$ws.Range("D41").Value = 24
$ws.Range("B42").Value = "global warming"
$ws.Range("C42").Value = "renewable energy"
$ws.Range("B43").Value = "global warming"
$ws.Range("C43").Value = "misinformation"
$ws.Range("B44").Value = "melting glaciers"
$ws.Range("C44").Value = "renewable energy"
$ws.Range("D44").Value = 2
$ws.Range("B45").Value = "melting glaciers"
$ws.Range("C45").Value = "misinformation"
$ws.Range("B46").Value = "renewable energy"
$ws.Range("C46").Value = "misinformation"
